# Add a new row 80 to each of the 4 worksheets (MID_LFT_#1, MID_LFT_#2,
# MID_PLT_#1, MID_PLT_#2), matching the new-row data appended in the
# source diff. Column A keeps the existing date/time number format used
# by the rest of the column; columns B-E are short hex-byte strings
# stored as text; columns F-I are plain numbers.

$wb = $excel.ActiveWorkbook

$rowData = @{
    1 = @{
        A = 45866.46563657407
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x34"
        E = "0x07"
        F = 400
        G = "5.68631262647113e+23"
        H = 308
        I = 7
    }
    2 = @{
        A = 45866.46563657407
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x34"
        E = "0x19"
        F = 380
        G = "5.68432987514711e+23"
        H = 308
        I = 25
    }
    3 = @{
        A = 45866.46563657407
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x60"
        E = "0x15"
        F = 110
        G = "5.68631262647113e+23"
        H = 96
        I = 15
    }
    4 = @{
        A = 45866.46563657407
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x77"
        E = "0x9"
        F = 130
        G = "5.68631262647113e+23"
        H = 119
        I = 9
    }
}

foreach ($sheetIndex in 1..4) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $data = $rowData[$sheetIndex]
    $newRow = 80

    $ws.Cells.Item($newRow, 1).Value = $data.A
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E

    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = [double]$data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}

Write-Host "Row 80 added to all 4 sheets"
